$wb = $excel.ActiveWorkbook
$wsLine = $wb.Worksheets.Item("Line")
$wsJumper = $wb.Worksheets.Item("Jumper")

# A new Line entry ("Line_12", connecting buses 9-10, disabled: u=0) is being added
# to the Line sheet as row 13, representing the jumper J9-10 as a de-energized line.
# Shift existing rows 13..20 down to 14..21 to make room, copying bottom-up so no
# data is overwritten before it's copied.
for ($r = 20; $r -ge 13; $r--) {
    $src = $wsLine.Range("A" + $r + ":U" + $r)
    $dst = $wsLine.Range("A" + ($r + 1) + ":U" + ($r + 1))
    $src.Copy($dst)
}

# Populate the newly-vacated row 13 with the new Line_12 record.
$wsLine.Cells.Item(13, 1).Value = 11          # uid
$wsLine.Cells.Item(13, 2).Value = "Line_12"   # idx
$wsLine.Cells.Item(13, 3).Value = 0           # u
$wsLine.Cells.Item(13, 4).Value = "Line_12"   # name
$wsLine.Cells.Item(13, 5).Value = 9           # bus1
$wsLine.Cells.Item(13, 6).Value = 10          # bus2
$wsLine.Cells.Item(13, 7).Value = 100         # Sn
$wsLine.Cells.Item(13, 8).Value = 60          # fn
$wsLine.Cells.Item(13, 9).Value = 138         # Vn1
$wsLine.Cells.Item(13, 10).Value = 138        # Vn2
$wsLine.Cells.Item(13, 11).Value = 0.03181    # r
$wsLine.Cells.Item(13, 12).Value = 0.0845     # x
$wsLine.Cells.Item(13, 13).Value = 0          # b
$wsLine.Cells.Item(13, 14).Value = 0          # g
$wsLine.Cells.Item(13, 15).Value = 0          # b1
$wsLine.Cells.Item(13, 16).Value = 0          # g1
$wsLine.Cells.Item(13, 17).Value = 0          # b2
$wsLine.Cells.Item(13, 18).Value = 0          # g2
$wsLine.Cells.Item(13, 19).Value = 0          # trans
$wsLine.Cells.Item(13, 20).Value = 1          # tap
$wsLine.Cells.Item(13, 21).Value = 0          # phi

# Update selections/active-tab: the Jumper tab was selected before; now the Line
# tab (with the newly-added row) becomes the active/selected tab.
[void]$wsJumper.Range("C3").Select()
[void]$wsLine.Activate()
[void]$wsLine.Range("C14").Select()
